$wb = $excel.ActiveWorkbook

# Activate the "Repayment schedule" worksheet (this becomes the active/selected tab)
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new blank column before column N (14th column), shifting
# the existing "Late" / "Outstanding"(heading) / "Outstanding" columns right.
# (Capture column M's width first, since Excel copies the left-neighbour's
# formatting/width into a freshly inserted column.)
$mWidth = $ws.Columns.Item(13).ColumnWidth
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = $mWidth

# Set the selection to K16 as the final cursor position on this sheet
$ws.Range("K16").Select()
